# Apply updated TPM values to sheet1, and delete the now-obsolete rows 6-9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (columns A through T).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$data = @(
    @("ECs",           "C3", "Itgam", "Resolving-Mac", 3, 1, 0.555934,           1.667802,    0.005745252779589096, 0.005745252779589094, 3, 1, 35.68243999999999, 107.04732, 1, 1, 19.83708159896,    178.53373439064,    0.005745252779589096, 0.005745252779589094),
    @("FAPs",          "C3", "Itgam", "Resolving-Mac", 3, 1, 79.68771233333334, 239.063137,  0.823525905561055,    0.823525905561055,    3, 1, 35.68243999999999, 107.04732, 1, 1, 2843.452014071426, 25591.06812664284,  0.823525905561055,    0.823525905561055),
    @("MuSCs",         "C3", "Itgam", "Resolving-Mac", 3, 1, 0.3446996666666666,1.034099,    0.003562269474506148, 0.003562269474506148, 3, 1, 35.68243999999999, 107.04732, 1, 1, 12.29972517385333, 110.69752656468,    0.003562269474506148, 0.003562269474506148),
    @("Resolving-Mac", "C3", "Itgam", "Resolving-Mac", 3, 1, 16.17571666666667, 48.52715,    0.1671665721848498,   0.1671665721848498,   3, 1, 35.68243999999999, 107.04732, 1, 1, 577.1890394153332, 5194.701354737999,  0.1671665721848498,   0.1671665721848498)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# Remove the old rows 6-9, which no longer exist in the updated data set.
$ws.Range("A6:T9").EntireRow.Delete()
